# Append the GymWorkouts entries recorded 02-08 April 2018 (commit:
# "Updating GymWorkouts file 08/04/2018") to the WeightTraining sheet.
# Columns: ExerciseId | DateId | Exercise Date | Exercise Month | Exercise Year |
#          Exercise Day | Exercise Name | Weight | Sets | Reps | TrainingArea
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ("730","85","4/2/2018","April","2018","Monday","Bench Press","85","5","5","Chest"),
    ("731","85","4/2/2018","April","2018","Monday","Overhead Press","50","5","5","Shoulders"),
    ("732","85","4/2/2018","April","2018","Monday","Barbell Row","90","5","5","Back"),
    ("733","85","4/2/2018","April","2018","Monday","Pull-Ups","93","5","5","Shoulders"),
    ("734","85","4/2/2018","April","2018","Monday","Bicep Curl","30","4","12","Arms"),
    ("735","85","4/2/2018","April","2018","Monday","Upright Rows","30","4","12","Arms"),
    ("736","86","4/4/2018","April","2018","Wednesday","Deadlift","125","5","5","Legs"),
    ("737","86","4/4/2018","April","2018","Wednesday","Hip adduction","70","4","8","Legs"),
    ("738","86","4/4/2018","April","2018","Wednesday","Hip abduction","70","4","8","Legs"),
    ("739","86","4/4/2018","April","2018","Wednesday","Leg Press","150","4","8","Legs"),
    ("740","87","4/5/2018","April","2018","Thursday","Pec Fly","120","5","5","Chest"),
    ("741","87","4/5/2018","April","2018","Thursday","Tricep Pull down","42.5","4","8","Arms"),
    ("742","87","4/5/2018","April","2018","Thursday","Hammer Curl","17.5","4","8","Arms"),
    ("743","87","4/5/2018","April","2018","Thursday","Seated Row","80","4","8","Back"),
    ("744","87","4/5/2018","April","2018","Thursday","Plank","0","4","30","Core"),
    ("745","87","4/5/2018","April","2018","Thursday","Left Plank","0","4","30","Core"),
    ("746","87","4/5/2018","April","2018","Thursday","Right Plank","0","4","30","Core"),
    ("747","88","4/7/2018","April","2018","Saturday","Incline Bench","75","5","5","Chest"),
    ("748","88","4/7/2018","April","2018","Saturday","Pec Fly","120","5","5","Chest"),
    ("749","88","4/7/2018","April","2018","Saturday","Shoulder Press","30","5","5","Shoulders"),
    ("750","88","4/7/2018","April","2018","Saturday","Shoulder Shrug","30","5","5","Shoulders"),
    ("751","88","4/7/2018","April","2018","Saturday","Heel-taps","0","4","12","Core"),
    ("752","88","4/7/2018","April","2018","Saturday","Raised leg circles","0","4","10","Core"),
    ("753","88","4/7/2018","April","2018","Saturday","Scissors","0","4","12","Core"),
    ("754","88","4/7/2018","April","2018","Saturday","Knee-Pull ins","0","4","10","Core"),
    ("755","88","4/7/2018","April","2018","Saturday","Flitter Kicks","0","4","20","Core"),
    ("756","89","4/8/2018","April","2018","Sunday","Bench Press","85","5","5","Chest"),
    ("757","89","4/8/2018","April","2018","Sunday","Overhead Press","50","5","5","Shoulders"),
    ("758","89","4/8/2018","April","2018","Sunday","Barbell Row","95","5","5","Back"),
    ("759","89","4/8/2018","April","2018","Sunday","Pull-Ups","93","5","5","Shoulders"),
    ("760","89","4/8/2018","April","2018","Sunday","Upright Rows","30","4","12","Shoulders"),
    ("761","89","4/8/2018","April","2018","Sunday","Bicep Curl","35","4","8","Arms")
)

$data = New-Object 'object[,]' $rows.Count,11
for ($r = 0; $r -lt $rows.Count; $r++) {
    for ($c = 0; $c -lt 11; $c++) {
        $data[$r,$c] = $rows[$r][$c]
    }
}

$firstRow = 731
$lastRow = $firstRow + $rows.Count - 1
$ws.Range("A$firstRow`:K$lastRow").Value = $data

# Mirror the author's final cursor position after the edit.
$ws.Range("C764").Select()

